$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("requiremens")
$ws2 = $wb.Worksheets.Item("Bugs")

# New rows on sheet1 ("requiremens").
# Shared-string indices must be minted in this order (14,15,16,17) to match
# the target: row9 A9 first, then row11 A11, then row10 A10, then sheet2 row4 A4.
$ws1.Cells.Item(9, 1).Value = "Cerebro - Need utility to create page dynamically with entries i give"
$ws1.Cells.Item(11, 1).Value = "Update Cerebro with all core java example"
$ws1.Cells.Item(11, 2).Value = "Completed"
$ws1.Cells.Item(10, 1).Value = "Update Cerebro with core java definitions"
$ws1.Cells.Item(10, 2).Value = "Completed"

# New row on sheet2 ("Bugs").
$ws2.Cells.Item(4, 1).Value = "Add Defintions web services -> service factory -> util class has 2 different method. But 2 methods are not needed, make it 1"

# Selection / active-sheet bookkeeping: requiremens keeps a stale selection at
# B12, then Bugs becomes the active (tabSelected) sheet with selection A5.
$ws1.Activate() | Out-Null
$ws1.Range("B12").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("A5").Select() | Out-Null
